# Commit: "Nguyen Trong Tien add link"
# Adds a new student row (Nguyen Trong Tien / 6690007) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row
$ws.Range("B3").Value = "Nguyen Trong Tien"
$ws.Range("C3").Value = 6690007

# Widen column B so the new name fits (matches the author's manual resize)
$ws.Columns.Item(2).ColumnWidth = 17

# Leave the selection on the last edited cell, like the saved workbook shows
$ws.Range("D3").Select()
